# Removing less than USD 5 price from extrapolation calibration because it is just a noise
# This updates the derived calibration columns (D:H) for the affected rows to reflect
# recalibration after excluding a sub-$5 price point considered noise.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row, D, E, F, G, H
$updates = @(
    @(2,  107476.9161419617,  0.003617399874283669, 0.173625874445852,  -1.589981041063056, 17.4178780006375),
    @(3,  108232.0758374846, -0.004542212730375906, 0.2037658022930758, -1.825199076932311, 19.90438322351699),
    @(4,  108837.6930382102, -0.0113524103336068,   0.2238400813988149, -1.665945255668843, 17.22781570930497),
    @(5,  109477.0727525339, -0.02092515166765615,  0.240629380352568,  -1.368125743226005, 12.0913129524412),
    @(7,  110603.2708207846, -0.03257082180835618,  0.2658790274177847, -1.599152670106467, 13.97573728883357),
    @(8,  112064.1073413354, -0.04930393833406367,  0.2283365855568399, -1.031343685244,    7.808993686208689),
    @(9,  113811.4031281088, -0.09629787212868798,  0.3892634046940351, -1.560341121503166, 9.3449970508922),
    @(10, 115121.562192494,  -0.1378937230725891,   0.4933568717092279, -1.666786479785098, 7.852648229366139),
    @(11, 117152.6102817787, -0.1493872439347294,   0.4902954498040062, -1.489340449343088, 6.808393445811004),
    @(16, 107486.8139978289,  0.04917852490688815,  0.1312076037811604, -0.4960488970660601, 5.763117108828423),
    @(17, 107389.6115682416,  0.03801546998358998,  0.1271071697411284, -0.8200296653877858, 5.875091313900017),
    @(20, 108399.3838202697,  0.0250789756816017,   0.1639977452644685, -0.3773779946812458, 4.89018952680705)
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 4).Value = $u[1]   # D
    $ws.Cells.Item($row, 5).Value = $u[2]   # E
    $ws.Cells.Item($row, 6).Value = $u[3]   # F
    $ws.Cells.Item($row, 7).Value = $u[4]   # G
    $ws.Cells.Item($row, 8).Value = $u[5]   # H
}
